{"js": "// Apply the Spanish translations described by the diff:\n//   \"GPU NVIDIA miner\"                       -> \"Minero de NVIDIA GPU\"\n//   \"GPU AMD miner\"                          -> \"Minero de GPU AMD\"\n//   \"For mining support please join\"         -> \"Para obtener soporte minero, \u00fanase\"\n//   \"\\u00a0EXCHANGES\" (nbsp + EXCHANGES)     -> \"\\u00a0INTERCAMBIOS\"\n//   \"Turn altcoins into Smartcash instantly\" -> \"Convertir altcoins en Smartcash al instante\"\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\"GPU NVIDIA miner\", \"Minero de NVIDIA GPU\"],\n  [\"GPU AMD miner\", \"Minero de GPU AMD\"],\n  [\"For mining support please join\", \"Para obtener soporte minero, \u00fanase\"],\n  [\"\\u00a0EXCHANGES\", \"\\u00a0INTERCAMBIOS\"],\n  [\"Turn altcoins into Smartcash instantly\", \"Convertir altcoins en Smartcash al instante\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Spanish translations described by the diff:\n#   \"GPU NVIDIA miner\"                       -> \"Minero de NVIDIA GPU\"\n#   \"GPU AMD miner\"                          -> \"Minero de GPU AMD\"\n#   \"For mining support please join\"         -> \"Para obtener soporte minero, \u00fanase\"\n#   \"<nbsp>EXCHANGES\"                        -> \"<nbsp>INTERCAMBIOS\"\n#   \"Turn altcoins into Smartcash instantly\" -> \"Convertir altcoins en Smartcash al instante\"\n\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\n$replacements = @(\n    @(\"GPU NVIDIA miner\", \"Minero de NVIDIA GPU\"),\n    @(\"GPU AMD miner\", \"Minero de GPU AMD\"),\n    @(\"For mining support please join\", \"Para obtener soporte minero, \u00fanase\"),\n    @($nbsp + \"EXCHANGES\", $nbsp + \"INTERCAMBIOS\"),\n    @(\"Turn altcoins into Smartcash instantly\", \"Convertir altcoins en Smartcash al instante\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $oldText\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Text = $newText\n    }\n}\n"}
